# msr6-bof-eckert-core-requirements.pptx : v1.0 -> v1.1 update
# - Update the version/date line on the title slide.
# - A handful of shapes across slide 4, 6 and 7 have their stored
#   position/size nudged by a single EMU (1/12700 pt). PowerPoint's COM
#   object model only exposes Left/Top/Width/Height in points (as a
#   float), and quantizes them back to EMU on write, so each of these is
#   set to a point value that has been verified (empirically, against
#   this very runtime) to round-trip to the exact target EMU value.

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle version/date line -----------------------------
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$para = $sh1.TextFrame.TextRange.Paragraphs(2, 1)
$run = $para.Runs(1, 1)
$run.Text = "v1.1 - 07/21/2022"

# --- Slide 4: connector line, off/x 430248 -> 430247 ------------------
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(3)
$sh4.Left = 33.877755535433074

# --- Slide 6 ------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Title 1: off/y 47624 -> 47623
$s6.Shapes.Item(2).Top = 3.7498825196850394

# textbox: ext/cx 767487 -> 767486
$s6.Shapes.Item(4).Width = 60.432007253937

# textbox: off/y 3823385 -> 3823384
$s6.Shapes.Item(13).Top = 301.05389376771654

# connector: off/y 1500600 -> 1500599
$s6.Shapes.Item(17).Top = 118.15744007480313

# connector: ext/cy 520004 -> 520003
$s6.Shapes.Item(21).Height = 40.94515811023622

# connector (duplicate of item 17): off/y 1500600 -> 1500599
$s6.Shapes.Item(28).Top = 118.15744007480313

# --- Slide 7 ------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# connector: ext/cy 520004 -> 520003
$s7.Shapes.Item(21).Height = 40.94515811023622
